$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated odds/scoreline values for 2025-02-04 FlashScore weekly games export.
# Each assignment corresponds to one changed cell from the source diff.

# Row 2
$ws.Range("G2").Value = 3.8
$ws.Range("I2").Value = 2.1
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 2.03
$ws.Range("AK2").Value = 19
# Row 4
$ws.Range("G4").Value = 1.48
$ws.Range("H4").Value = 3.8
$ws.Range("AA4").Value = 2.5
$ws.Range("AB4").Value = 1.5
$ws.Range("AF4").Value = 9.5
$ws.Range("AK4").Value = 26
$ws.Range("AQ4").Value = 101
# Row 5
$ws.Range("G5").Value = 2.63
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("S5").Value = 2.3
$ws.Range("T5").Value = 1.62
$ws.Range("W5").Value = 4.33
$ws.Range("X5").Value = 1.22
$ws.Range("Y5").Value = 1.5
$ws.Range("Z5").Value = 2.5
$ws.Range("AA5").Value = 1.91
$ws.Range("AB5").Value = 1.8
$ws.Range("AC5").Value = 7.5
$ws.Range("AH5").Value = 34
# Row 8
$ws.Range("G8").Value = 3.55
$ws.Range("H8").Value = 3.4
$ws.Range("K8").Value = 2.12
$ws.Range("L8").Value = 2.57
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 7.1
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 3.1
$ws.Range("S8").Value = 1.93
$ws.Range("T8").Value = 1.78
$ws.Range("W8").Value = 3.2
$ws.Range("X8").Value = 1.3
$ws.Range("Y8").Value = 1.4
$ws.Range("Z8").Value = 2.75
$ws.Range("AA8").Value = 1.8
$ws.Range("AB8").Value = 1.91
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 19
$ws.Range("AE8").Value = 12
$ws.Range("AG8").Value = 32
$ws.Range("AH8").Value = 40
$ws.Range("AI8").Value = 7.1
$ws.Range("AJ8").Value = 6.5
$ws.Range("AK8").Value = 14.5
$ws.Range("AL8").Value = 70
$ws.Range("AM8").Value = 600
$ws.Range("AN8").Value = 7.1
$ws.Range("AO8").Value = 9.25
$ws.Range("AP8").Value = 8.5
$ws.Range("AQ8").Value = 17.5
$ws.Range("AR8").Value = 16
$ws.Range("AS8").Value = 28
